$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3003
$ws.Range("F3").Value = 6430
$ws.Range("F4").Value = 2549
$ws.Range("F6").Value = 652
$ws.Range("F7").Value = 85
$ws.Range("F9").Value = 3144
$ws.Range("F10").Value = 366
$ws.Range("F12").Value = 7680
$ws.Range("F13").Value = 379
$ws.Range("F18").Value = 19
$ws.Range("F19").Value = 486
$ws.Range("F20").Value = 9433
$ws.Range("F21").Value = 25
$ws.Range("F23").Value = 271
$ws.Range("F26").Value = 29
$ws.Range("F27").Value = 128
$ws.Range("F28").Value = 130
$ws.Range("F31").Value = 74
$ws.Range("F33").Value = 2623
$ws.Range("F36").Value = 2051
$ws.Range("F39").Value = 3962
$ws.Range("F41").Value = 47
$ws.Range("F43").Value = 111
$ws.Range("F44").Value = 254
$ws.Range("F45").Value = 60
$ws.Range("F46").Value = 18
$ws.Range("F47").Value = 69
$ws.Range("F48").Value = 43
$ws.Range("F49").Value = 65
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 54
$ws.Range("F8").Value = 158
$ws.Range("F10").Value = 5
$ws.Range("F12").Value = 22
$ws.Range("F13").Value = 4
$ws.Range("F16").Value = 12
$ws.Range("F21").Value = 13
$ws.Range("F23").Value = 11
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 54
$ws.Range("F3").Value = 3004
$ws.Range("F6").Value = 6430
$ws.Range("F7").Value = 2549
$ws.Range("F8").Value = 158
$ws.Range("F10").Value = 652
$ws.Range("F11").Value = 85
$ws.Range("F13").Value = 3144
$ws.Range("F14").Value = 366
$ws.Range("F18").Value = 7680
$ws.Range("F19").Value = 379
$ws.Range("F23").Value = 19
$ws.Range("F24").Value = 9434
$ws.Range("F25").Value = 25
$ws.Range("F26").Value = 271
$ws.Range("F28").Value = 29
$ws.Range("F29").Value = 128
$ws.Range("F30").Value = 130
$ws.Range("F33").Value = 74
$ws.Range("F35").Value = 2623
$ws.Range("F36").Value = 2051
$ws.Range("F40").Value = 3962
$ws.Range("F42").Value = 47
$ws.Range("F44").Value = 111
$ws.Range("F45").Value = 254
$ws.Range("F46").Value = 60
$ws.Range("F47").Value = 69
$ws.Range("F48").Value = 43
$ws.Range("F49").Value = 65
